{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst OLD_HYBRID = \"Inklusions- und Lernstandszeugnisse (Hybrid-Berichte)\";\nconst NEW_HYBRID = \"Inklusions- und Lernstandszeugnisse \";\n\nconst OLD_AGGT = \"Die F\u00e4chergruppe AGGT wurde hinzugef\u00fcgt\";\nconst NEW_AGGT = \"Die F\u00e4chergruppen AGGT  und ZUV wurden hinzugef\u00fcgt.\";\n\nfor (const p of paragraphs.items) {\n  if (p.text === OLD_HYBRID) {\n    p.insertText(NEW_HYBRID, \"Replace\");\n  } else if (p.text === OLD_AGGT) {\n    p.insertText(NEW_AGGT, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# \"Inklusions- und Lernstandszeugnisse (Hybrid-Berichte)\" -> \"Inklusions- und Lernstandszeugnisse \"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"Inklusions- und Lernstandszeugnisse (Hybrid-Berichte)\", $false, $false, $false, $false, $false, $true, 1, $false, \"Inklusions- und Lernstandszeugnisse \", 2)\n\n# \"Die F\u00e4chergruppe AGGT wurde hinzugef\u00fcgt\" -> \"Die F\u00e4chergruppen AGGT  und ZUV wurden hinzugef\u00fcgt.\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"Die F\u00e4chergruppe AGGT wurde hinzugef\u00fcgt\", $false, $false, $false, $false, $false, $true, 1, $false, \"Die F\u00e4chergruppen AGGT  und ZUV wurden hinzugef\u00fcgt.\", 2)\n"}
